# This script applies the "stock-screener" refresh described by the commit:
# the values of the symbol columns (B, C, E, F) are updated for rows 2-37,
# and 11 new rows (38-48) are appended with additional "support Zone" (C) entries,
# continuing the running index in column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letter -> 1-based index, matching the sheet's A:F layout.
$colIndex = @{ A = 1; B = 2; C = 3; D = 4; E = 5; F = 6 }

# Cell A2 already carries the bold / centered / thin-bordered style used for every
# value in the running-index column (A) and the header row; reuse it (via
# copy/PasteSpecial of formats only) so newly appended index cells (A38:A48) pick
# up the same style instead of creating new, duplicate style entries.
$styleSource = $ws.Cells.Item(2, 1)

# Each entry below lists only the columns whose value actually changes for that
# row; "New = $true" marks a brand-new row that needs the index-column style
# applied after its value is written.
$rowEdits = @(
    @{ Row = 2; B = "NSE:COROMANDEL"; C = "NSE:ADANIGREEN"; E = "NSE:ADANIPORTS"; F = "" },
    @{ Row = 3; B = ""; C = "NSE:ADANIPORTS"; E = "NSE:AXISBANK"; F = "" },
    @{ Row = 4; B = ""; C = "NSE:AEROFLEX"; E = "NSE:CANFINHOME" },
    @{ Row = 5; B = ""; C = "NSE:AKZOINDIA"; E = "NSE:DALBHARAT" },
    @{ Row = 6; B = ""; C = "NSE:ANUP"; E = "NSE:ICICIPRULI" },
    @{ Row = 7; B = ""; C = "NSE:APCL"; E = "NSE:NESTLEIND" },
    @{ Row = 8; B = ""; C = "NSE:ARTEMISMED"; E = "NSE:PETRONET" },
    @{ Row = 9; B = ""; C = "NSE:BEARDSELL"; E = "NSE:RELIANCE" },
    @{ Row = 10; B = ""; C = "NSE:CAMLINFINE"; E = "" },
    @{ Row = 11; C = "NSE:CAMS"; E = "" },
    @{ Row = 12; C = "NSE:CENTEXT"; E = "" },
    @{ Row = 13; C = "NSE:CHEVIOT"; E = "" },
    @{ Row = 14; C = "NSE:DBCORP"; E = "" },
    @{ Row = 15; C = "NSE:DBSTOCKBRO"; E = "" },
    @{ Row = 16; C = "NSE:DEN"; E = "" },
    @{ Row = 17; C = "NSE:DHANI"; E = "" },
    @{ Row = 18; C = "NSE:DSSL"; E = "" },
    @{ Row = 19; C = "NSE:EMMBI"; E = "" },
    @{ Row = 20; C = "NSE:GEOJITFSL"; E = "" },
    @{ Row = 21; C = "NSE:GPIL" },
    @{ Row = 22; C = "NSE:GRAPHITE" },
    @{ Row = 23; C = "NSE:HATHWAY" },
    @{ Row = 24; C = "NSE:HEIDELBERG" },
    @{ Row = 25; C = "NSE:HINDZINC" },
    @{ Row = 26; C = "NSE:HITECHCORP" },
    @{ Row = 27; C = "NSE:INDSWFTLAB" },
    @{ Row = 28; C = "NSE:INDSWFTLTD" },
    @{ Row = 29; C = "NSE:IOLCP" },
    @{ Row = 30; C = "NSE:JKTYRE" },
    @{ Row = 31; C = "NSE:KANPRPLA" },
    @{ Row = 32; C = "NSE:KELLTONTEC" },
    @{ Row = 33; C = "NSE:KUANTUM" },
    @{ Row = 34; C = "NSE:MAHAPEXLTD" },
    @{ Row = 35; C = "NSE:MAPMYINDIA" },
    @{ Row = 36; C = "NSE:MATRIMONY" },
    @{ Row = 37; C = "NSE:MOKSH" },
    @{ Row = 38; New = $true; A = 36; C = "NSE:MSUMI" },
    @{ Row = 39; New = $true; A = 37; C = "NSE:MURUDCERA" },
    @{ Row = 40; New = $true; A = 38; C = "NSE:NETWORK18" },
    @{ Row = 41; New = $true; A = 39; C = "NSE:NLCINDIA" },
    @{ Row = 42; New = $true; A = 40; C = "NSE:NYKAA" },
    @{ Row = 43; New = $true; A = 41; C = "NSE:ORIENTPPR" },
    @{ Row = 44; New = $true; A = 42; C = "NSE:PFS" },
    @{ Row = 45; New = $true; A = 43; C = "NSE:POLYCAB" },
    @{ Row = 46; New = $true; A = 44; C = "NSE:PPAP" },
    @{ Row = 47; New = $true; A = 45; C = "NSE:RADIOCITY" },
    @{ Row = 48; New = $true; A = 46; C = "NSE:RELIGARE" }
)

foreach ($edit in $rowEdits) {
    foreach ($col in $colIndex.Keys) {
        if ($edit.ContainsKey($col)) {
            $ws.Cells.Item($edit.Row, $colIndex[$col]).Value = $edit[$col]
        }
    }
    if ($edit.ContainsKey('New') -and $edit.New) {
        $styleSource.Copy()
        $ws.Cells.Item($edit.Row, 1).PasteSpecial(-4122)
    }
}
